# Apply the "Penalty Reward System" forecast refresh:
#  - Shift every Week_Start_Date in the "Forecast Comparison" sheet forward by
#    one week (row 2 gets what used to be row 3's date, etc., with a new final
#    week appended).
#  - Replace the MyForecast (column D) numbers with the refreshed forecast.
#  - Update the "Summary" sheet metrics so they reflect the new forecast data.
#
# NOTE: the Week_Start_Date / date-ish Summary values and the numeric-looking
# Summary strings are stored as plain TEXT in the workbook (not real dates or
# numbers). Assigning a bare date-like or number-like string via .Value makes
# Excel auto-convert it to a real date serial / number, so each such value is
# written with a leading single-quote (the standard "force text" prefix) to
# keep it a text value, exactly like the source data.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New Week_Start_Date values (each row takes on what used to be the following
# row's date; the series extends one additional week at the end).
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast values.
$newForecast = @(174, 234, 180, 181, 187, 186, 188, 197, 189, 194, 195, 190, 190, 182, 180, 178)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 2).Value = "'" + $newDates[$i]
    $wsForecast.Cells.Item($row, 4).Value = $newForecast[$i]
}

# Update the Summary sheet to reflect the refreshed forecast.
$wsSummary.Range("B2").Value  = "2023-01-22 to 2025-01-05"
$wsSummary.Range("B4").Value  = "'449"
$wsSummary.Range("B5").Value  = "'178"
$wsSummary.Range("B6").Value  = "'162"
$wsSummary.Range("B8").Value  = "16912 units"
$wsSummary.Range("B9").Value  = "'3026"
$wsSummary.Range("B10").Value = "'1528"
$wsSummary.Range("B11").Value = "'769"
$wsSummary.Range("B12").Value = "'234"
$wsSummary.Range("B13").Value = "'2025-01-19"
$wsSummary.Range("B14").Value = "'174"
$wsSummary.Range("B15").Value = "'2025-01-12"
